$wb = $excel.ActiveWorkbook

# The workbook reports localization status for the zh-cn / de-de handoff
# files. The "Overview" sheet mirrors the same status value (columns E
# and F, row 2) that the "zh-cn" / "de-de" detail sheets carry in their
# "Status" column (column C, row 2). Move the status along from
# "Ready for handoff" to "In Translation" everywhere it appears.
$newStatus = "In Translation"

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

$wsOverview.Range("E2").Value = $newStatus
$wsOverview.Range("F2").Value = $newStatus
$wsZhCn.Range("C2").Value = $newStatus
$wsDeDe.Range("C2").Value = $newStatus

# The status column on each sheet is sized to fit its contents; shrink
# the now-shorter columns back down to fit "In Translation". (Use
# numeric column indices with Columns.Item - letter indices are not
# reliable in this runtime.)
$wsOverview.Columns.Item(5).ColumnWidth = 12.5
$wsOverview.Columns.Item(6).ColumnWidth = 12.5
$wsZhCn.Columns.Item(3).ColumnWidth = 12.5
$wsDeDe.Columns.Item(3).ColumnWidth = 12.5

Write-Host "Updated status to 'In Translation' on Overview, zh-cn, de-de"
